# Auto-generated Excel COM-interop script to apply cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) values are numeric-looking strings that must remain exact text ---
# Force a text number format first so Excel does not auto-convert/round these strings,
# then restore the default 'Normal' style so the saved file matches the original formatting.
$priceCells = @{
    'D2' = '26.438.31'
    'D3' = '1.769.80'
    'D5' = '1.003'
    'D6' = '306.65'
    'D7' = '0.4290'
    'D9' = '0.07214'
    'D10' = '0.8487'
    'D12' = '1.768.08'
    'D13' = '6.426'
    'D14' = '5.231'
    'D15' = '0.06922'
    'D16' = '1.004'
    'D17' = '79.27'
    'D18' = '0.000008684'
    'D19' = '1.003'
    'D20' = '15.00'
    'D21' = '26.440.75'
    'D23' = '11.27'
    'D24' = '2.009.21'
    'D25' = '152.11'
    'D26' = '1.872'
    'D28' = '5.077'
    'D29' = '114.40'
    'D30' = '1.731'
    'D31' = '0.08968'
    'D32' = '0.7227'
    'D33' = '1.111'
    'D35' = '2.749'
    'D36' = '1.003'
    'D38' = '0.05161'
    'D39' = '0.01891'
    'D40' = '0.4921'
    'D41' = '0.1603'
    'D42' = '2.578'
    'D43' = '6.258'
    'D44' = '7.999'
    'D45' = '104.77'
    'D46' = '1.003'
    'D48' = '0.06198'
    'D49' = '0.4489'
    'D50' = '1.590'
    'D51' = '1.739'
}
foreach ($ref in $priceCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $priceCells[$ref]
    $cell.Style = "Normal"
}

# --- Other changed cells (Coin name, Link, Volume(1h)) ---
$otherCells = @{
    'E2' = '  -3.83%  '
    'E3' = '  -3.12%  '
    'E4' = '  +0.25%  '
    'E5' = '  +0.18%  '
    'E6' = '  -1.99%  '
    'E7' = '  +0.61%  '
    'E8' = '  +1.23%  '
    'E9' = '  -0.03%  '
    'E10' = '  -1.86%  '
    'E11' = '  -1.36%  '
    'E12' = '  -7.01%  '
    'E13' = '  -0.83%  '
    'E14' = '  -3.03%  '
    'E15' = '  -0.18%  '
    'E16' = '  +0.07%  '
    'E17' = '  -2.04%  '
    'E18' = '  -2.95%  '
    'E19' = '  +0.15%  '
    'E20' = '  -2.69%  '
    'E21' = '  -4.92%  '
    'E22' = '  -0.53%  '
    'E23' = '  +3.52%  '
    'E24' = '  -6.67%  '
    'B25' = 'Monero'
    'C25' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'E25' = '  -2.16%  '
    'B26' = 'Toncoin'
    'C26' = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
    'E26' = '  -5.89%  '
    'E27' = '  -3.36%  '
    'E28' = '  -1.26%  '
    'E29' = '  +0.08%  '
    'E30' = '  -3.61%  '
    'E31' = '  +0.74%  '
    'E32' = '  -3.53%  '
    'E33' = '  -1.00%  '
    'E34' = '  -4.98%  '
    'B35' = 'HuobiToken'
    'C35' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'E35' = '  -7.93%  '
    'B36' = 'Frax'
    'C36' = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
    'E36' = '  +0.18%  '
    'E37' = '  -0.68%  '
    'E38' = '  -2.07%  '
    'E39' = '  -1.96%  '
    'E40' = '  -3.18%  '
    'E42' = '  -7.72%  '
    'E43' = '  -2.91%  '
    'E44' = '  -4.35%  '
    'E45' = '  -1.65%  '
    'E46' = '  +0.23%  '
    'E47' = '  -3.18%  '
    'E48' = '  -4.09%  '
    'E49' = '  -4.33%  '
    'E50' = '  -1.63%  '
    'E51' = '  +2.63%  '
}
foreach ($ref in $otherCells.Keys) {
    $ws.Range($ref).Value = $otherCells[$ref]
}

Write-Host "Applied crypto data refresh."
